$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new blank row at row 7 (pushes jogging_km/jogging_min/mood/drinks/submit/note down by one)
$ws.Rows("7:7").Insert()

# New row 7: jogging_at_all multiple-choice question
$ws.Range("C7").Value = "mc"
$ws.Range("D7").Value = "jogging_at_all"
$ws.Range("F7").Value = "Did you go jogging today?"
$ws.Range("G7").Value = "yes"
$ws.Range("H7").Value = "no"
$ws.Range("J7").Clear()

# Header: rename "skipif" column to "showif"
$ws.Range("I1").Value = "showif"

# Row 8 (jogging_km, shifted down from 7): add showif condition, taller row
$ws.Range("I8").Value = "tail(jogging_at_all, 1) == 1"
$ws.Rows("8:8").RowHeight = 45

# Row 9 (jogging_min, shifted down from 8): add showif condition, taller row
$ws.Range("I9").Value = "tail(jogging_at_all, 1) == 1"
$ws.Rows("9:9").RowHeight = 45

# Row 10 (was an empty gap row, now holds the new feelgood question)
$ws.Range("C10").Value = "mc_button"
$ws.Range("D10").Value = "feelgood"
$ws.Range("F10").Value = "Did it feel good?"
$ws.Range("G10").Value = "yes"
$ws.Range("H10").Value = "no"
$ws.Range("I10").Value = "tail(jogging_at_all, 1) == 1"
$ws.Rows("10:10").RowHeight = 45

# Selection, as recorded in the saved workbook
$ws.Range("I10").Select()

$wb.Save()
